$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "users/*" endpoint rows (text unchanged, but shared-string
# indices shift because the "flashcards" / "flashcards/:id" strings were
# removed earlier in the table).
$ws.Range("C5").Value = "users/getdetails"
$ws.Range("C6").Value = "users/update/:id"
$ws.Range("C7").Value = "users/delete/:id"

# New, more specific Flashcard API endpoints.
$ws.Range("C16").Value = "flashcards/getAll"
$ws.Range("C17").Value = "flashcards/get/:id"
$ws.Range("C18").Value = "flashcards/:deckId"
$ws.Range("C19").Value = "flashcards/update/:deckId"
$ws.Range("C20").Value = "flashcards/delete/:deckid"

# Deck endpoints (text unchanged, indices shift).
$ws.Range("C21").Value = "decks/getAll"
$ws.Range("C22").Value = "decks/get/:id"
$ws.Range("C23").Value = "decks/post"
$ws.Range("C24").Value = "decks/update/:id"
$ws.Range("C25").Value = "decks/delete/:id"

# Restore the view to show the new Flashcard rows (as recorded in the
# workbook's sheetView after the edit).
$ws.Range("C20").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
